# Quarterly database update + "read_price" shift:
# every quarter-series on the sheet drops its oldest (leftmost, column E)
# quarter and gains one new quarter (rightmost, column N) — i.e. all values
# shift one column to the left and a freshly-read value lands in column N.
#
# This touches:
#   - the two header rows (8 and 24) that carry the quarter labels
#   - the nine data rows that carry one quarter-series each

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header rows: shift the quarter labels left, add the new quarter ---
$newQuarterLabel = "فصل چهارم منتهی به 1401/12"

foreach ($row in @(8, 24)) {
    $old = @()
    for ($col = 5; $col -le 14; $col++) {
        $old += $ws.Cells.Item($row, $col).Value2
    }
    for ($col = 5; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = $old[$col - 5 + 1]
    }
    $ws.Cells.Item($row, 14).Value = $newQuarterLabel
}

# --- 2. Data rows: shift each quarter-series left, append freshly read value ---
$newLastValue = @{
    10 = 485892
    13 = 86011
    15 = -2329
    16 = 15946
    17 = 173390
    19 = 322007
    20 = 1080917
    26 = 435
    27 = 311
}

foreach ($row in $newLastValue.Keys) {
    $old = @()
    for ($col = 5; $col -le 14; $col++) {
        $old += $ws.Cells.Item($row, $col).Value2
    }
    for ($col = 5; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = $old[$col - 5 + 1]
    }
    $ws.Cells.Item($row, 14).Value = $newLastValue[$row]
}
